# Natmi following Dr Hou advice
# Add FAPs/ECs cell-type rows to the Calca-Calcr LR-pair sheet and
# refresh the derived-specificity table to reflect the new cell-type set.
# (sCs, Calca, Calcr) row becomes 4 sender/receiver combinations:
#   FAPs->ECs, FAPs->sCs, sCs->ECs, sCs->sCs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Calca"
$ws.Cells.Item(2, 3).Value = "Calcr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.06082199999999999
$ws.Cells.Item(2, 8).Value = 0.182466
$ws.Cells.Item(2, 9).Value = 0.1419252066466042
$ws.Cells.Item(2, 10).Value = 0.1419252066466042
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.096149
$ws.Cells.Item(2, 14).Value = 0.288447
$ws.Cells.Item(2, 15).Value = 0.2085078145161115
$ws.Cells.Item(2, 16).Value = 0.2085078145161116
$ws.Cells.Item(2, 17).Value = 0.005847974478
$ws.Cells.Item(2, 18).Value = 0.052631770302
$ws.Cells.Item(2, 19).Value = 0.02959251466263094
$ws.Cells.Item(2, 20).Value = 0.02959251466263094

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Calca"
$ws.Cells.Item(3, 3).Value = "Calcr"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.06082199999999999
$ws.Cells.Item(3, 8).Value = 0.182466
$ws.Cells.Item(3, 9).Value = 0.1419252066466042
$ws.Cells.Item(3, 10).Value = 0.1419252066466042
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.36498
$ws.Cells.Item(3, 14).Value = 1.09494
$ws.Cells.Item(3, 15).Value = 0.7914921854838884
$ws.Cells.Item(3, 16).Value = 0.7914921854838886
$ws.Cells.Item(3, 17).Value = 0.02219881356
$ws.Cells.Item(3, 18).Value = 0.19978932204
$ws.Cells.Item(3, 19).Value = 0.1123326919839732
$ws.Cells.Item(3, 20).Value = 0.1123326919839732

# Row 4
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Calca"
$ws.Cells.Item(4, 3).Value = "Calcr"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3677276666666667
$ws.Cells.Item(4, 8).Value = 1.103183
$ws.Cells.Item(4, 9).Value = 0.8580747933533959
$ws.Cells.Item(4, 10).Value = 0.8580747933533959
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.096149
$ws.Cells.Item(4, 14).Value = 0.288447
$ws.Cells.Item(4, 15).Value = 0.2085078145161115
$ws.Cells.Item(4, 16).Value = 0.2085078145161116
$ws.Cells.Item(4, 17).Value = 0.03535664742233333
$ws.Cells.Item(4, 18).Value = 0.318209826801
$ws.Cells.Item(4, 19).Value = 0.1789152998534806
$ws.Cells.Item(4, 20).Value = 0.1789152998534806

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Calca"
$ws.Cells.Item(5, 3).Value = "Calcr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.3677276666666667
$ws.Cells.Item(5, 8).Value = 1.103183
$ws.Cells.Item(5, 9).Value = 0.8580747933533959
$ws.Cells.Item(5, 10).Value = 0.8580747933533959
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.36498
$ws.Cells.Item(5, 14).Value = 1.09494
$ws.Cells.Item(5, 15).Value = 0.7914921854838884
$ws.Cells.Item(5, 16).Value = 0.7914921854838886
$ws.Cells.Item(5, 17).Value = 0.13421324378
$ws.Cells.Item(5, 18).Value = 1.20791919402
$ws.Cells.Item(5, 19).Value = 0.6791594934999152
$ws.Cells.Item(5, 20).Value = 0.6791594934999153
